# Revert to 2.1.1 files
#
# 1) Add a new "Texas Data" worksheet between "IEA Data" and "HPEbP" with a
#    set of explanatory notes (styled with the accent-5 theme font color).
# 2) Fix the HPEbP!B3 efficiency formula: it was incorrectly including waste
#    heat (46) in the energy-balance denominator. Drop it so the formula
#    reads 118/(162+2) instead of 118/(162+2+46); all the cells to the right
#    (C3:AI3) reference B3 (directly or via the shared formula chain) so they
#    recalculate automatically.
# 3) Restore the various sheet selections / active tab to match the saved
#    state of the edited workbook.

$wb = $excel.ActiveWorkbook

# --- selection bookkeeping on the existing sheets -------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("B14").Select()

$wsIea = $wb.Worksheets.Item("IEA Data")
$wsIea.Range("E18").Select()

# --- insert the new "Texas Data" sheet right after "IEA Data" ------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Texas Data"
$newSheet.Move([Type]::Missing, $wb.Worksheets.Item("IEA Data"))

$ws = $wb.Worksheets.Item("Texas Data")

# Style the whole used range with the accent5-colored font used for the note.
$ws.Range("A1:I17").Font.ThemeColor = 9

$ws.Range("A1").Value = "There is no reason that these number should be different for Texas."
$ws.Range("A3").Value = "However, I did find an error in their calculations. "
$ws.Range("A5").Value = "They were included waste heat as an energy balance input."
$ws.Range("B6").Value = "for example, page 228 of the NREL report shows gas production as 162 kBtu gas + 2 kBtu electricity = 118 kBtu hydrogen + 46 kBtu waste heat"
$ws.Range("B7").Value = "so, the efficiency (output hydrogen energy vs input energy) would be 118/(162+2)=72%"
$ws.Range("B8").Value = "previously, this spreadhseet (cell 'HPEbP'B3) was calculating the efficiency as 118/(162+2+46)=56%"
$ws.Range("B10").Value = "the IEA number for natural gas reforming efficiency is 76%, so that's a good check that their initial calculation was wrong. "
$ws.Range("A12").Value = "Their other calculations did not include the same mistake."

$ws.Range("A13").Select()

# --- fix the HPEbP natural-gas-reforming efficiency formula ---------------
$wsHpebp = $wb.Worksheets.Item("HPEbP")
$wsHpebp.Range("B3").Formula = "=118/(162+2)"

# HPEbP ends up the active/visible tab in the saved workbook.
$wsHpebp.Activate()
$wsHpebp.Range("C12").Select()
